$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.049.63'
$ws.Range("E2").Value = '  +2.82%  '
$ws.Range("D3").Value = '3.801.63'
$ws.Range("E3").Value = '  +0.66%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''705.60'
$ws.Range("E5").Value = '  +11.94%  '
$ws.Range("D6").Value = '''172.69'
$ws.Range("E6").Value = '  +4.27%  '
$ws.Range("D7").Value = '3.800.41'
$ws.Range("E7").Value = '  +0.59%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '''0.526'
$ws.Range("E9").Value = '  +1.13%  '
$ws.Range("D10").Value = '''0.163'
$ws.Range("E10").Value = '  +3.14%  '
$ws.Range("D11").Value = '''7.48'
$ws.Range("E11").Value = '  +10.24%  '
$ws.Range("D12").Value = '''0.461'
$ws.Range("E12").Value = '  +1.01%  '
$ws.Range("E13").Value = '  +9.48%  '
$ws.Range("D14").Value = '''36.23'
$ws.Range("E14").Value = '  +3.87%  '
$ws.Range("D15").Value = '4.443.41'
$ws.Range("E15").Value = '  +0.89%  '
$ws.Range("D16").Value = '3.803.30'
$ws.Range("E16").Value = '  +0.34%  '
$ws.Range("D17").Value = '71.104.17'
$ws.Range("E17").Value = '  +3.05%  '
$ws.Range("D18").Value = '''17.90'
$ws.Range("E18").Value = '  +1.20%  '
$ws.Range("D19").Value = '''7.23'
$ws.Range("E19").Value = '  +3.07%  '
$ws.Range("E20").Value = '  +0.56%  '
$ws.Range("D21").Value = '''11.23'
$ws.Range("E21").Value = '  +17.95%  '
$ws.Range("D22").Value = '''484.46'
$ws.Range("E22").Value = '  +3.50%  '
$ws.Range("D23").Value = '''0.715'
$ws.Range("E23").Value = '  +1.80%  '
$ws.Range("D24").Value = '''0.0000147'
$ws.Range("E24").Value = '  +3.32%  '
$ws.Range("D25").Value = '''83.66'
$ws.Range("E25").Value = '  +1.90%  '
$ws.Range("D26").Value = '''12.38'
$ws.Range("E26").Value = '  +1.89%  '
$ws.Range("D27").Value = '''10.56'
$ws.Range("E27").Value = '  +4.19%  '
$ws.Range("D28").Value = '''2.18'
$ws.Range("E28").Value = '  +2.63%  '
$ws.Range("D29").Value = '3.954.76'
$ws.Range("E29").Value = '  +0.77%  '
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("D31").Value = '''3.13'
$ws.Range("E31").Value = '  +17.23%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '''2.30'
$ws.Range("E32").Value = '  +0.47%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").Value = '''7.56'
$ws.Range("E33").Value = '  +6.16%  '
$ws.Range("D34").Value = '''29.60'
$ws.Range("E34").Value = '  +3.99%  '
$ws.Range("E35").Value = '  +0.30%  '
$ws.Range("D36").Value = '''9.24'
$ws.Range("E36").Value = '  +3.85%  '
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("D38").Value = '3.754.19'
$ws.Range("E38").Value = '  +0.75%  '
$ws.Range("D39").Value = '''0.103'
$ws.Range("E39").Value = '  +2.04%  '
$ws.Range("D40").Value = '''3.50'
$ws.Range("E40").Value = '  +6.98%  '
$ws.Range("D41").Value = '''5.98'
$ws.Range("E41").Value = '  +3.17%  '
$ws.Range("D42").Value = '''2.22'
$ws.Range("E42").Value = '  +11.06%  '
$ws.Range("D43").Value = '''0.000329'
$ws.Range("E43").Value = '  +24.12%  '
$ws.Range("D44").Value = '''0.969'
$ws.Range("E44").Value = '  +0.69%  '
$ws.Range("E45").Value = '  +0.15%  '
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").Value = '''161.63'
$ws.Range("E47").Value = '  +3.49%  '
$ws.Range("D48").Value = '''49.42'
$ws.Range("E48").Value = '  +5.11%  '
$ws.Range("D49").Value = '''45.13'
$ws.Range("E49").Value = '  +3.34%  '
$ws.Range("D50").Value = '''1.40'
$ws.Range("E50").Value = '  -0.56%  '
$ws.Range("E51").Value = '  +2.41%  '
